$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(40, 8).Value = 1881.8182
$ws.Cells.Item(40, 9).Value = 1712.5
$ws.Cells.Item(40, 11).Value = 1712.5
$ws.Cells.Item(40, 13).Value = -1537.5
$ws.Cells.Item(43, 8).Value = 1190.1
$ws.Cells.Item(43, 9).Value = 750
$ws.Cells.Item(43, 10).Value = 1483.5
$ws.Cells.Item(43, 11).Value = 750
$ws.Cells.Item(43, 12).Value = 1483.5
$ws.Cells.Item(43, 13).Value = -681
$ws.Cells.Item(43, 14).Value = -1621.5
$ws.Cells.Item(62, 8).Value = 1955.8
$ws.Cells.Item(62, 9).Value = 1953.1111
$ws.Cells.Item(62, 10).Value = 1980
$ws.Cells.Item(62, 11).Value = 1953.1111
$ws.Cells.Item(62, 12).Value = 1980
$ws.Cells.Item(62, 13).Value = -1329.1111
$ws.Cells.Item(62, 14).Value = -3228
$ws.Cells.Item(65, 8).Value = 1955.8
$ws.Cells.Item(65, 9).Value = 1953.1111
$ws.Cells.Item(65, 10).Value = 1980
$ws.Cells.Item(65, 11).Value = 9765.5555
$ws.Cells.Item(65, 12).Value = 9900
$ws.Cells.Item(65, 13).Value = -6645.5555
$ws.Cells.Item(65, 14).Value = -16140
$ws.Cells.Item(92, 8).Value = 634.8823
$ws.Cells.Item(92, 9).Value = 549.5625
$ws.Cells.Item(92, 10).Value = 2000
$ws.Cells.Item(92, 11).Value = 549.5625
$ws.Cells.Item(92, 12).Value = 2000
$ws.Cells.Item(92, 13).Value = 698.4375
$ws.Cells.Item(92, 14).Value = -4496
$ws.Cells.Item(106, 8).Value = 5199.8887
$ws.Cells.Item(106, 9).Value = 5199.8887
$ws.Cells.Item(106, 11).Value = 5199.8887
$ws.Cells.Item(106, 13).Value = -4568.8887
$ws.Cells.Item(116, 8).Value = 2438.6667
$ws.Cells.Item(116, 9).Value = 2002
$ws.Cells.Item(116, 10).Value = 2657
$ws.Cells.Item(116, 11).Value = 2002
$ws.Cells.Item(116, 12).Value = 2657
$ws.Cells.Item(116, 13).Value = 1440
$ws.Cells.Item(116, 14).Value = -9541
$ws.Cells.Item(137, 8).Value = 2086.6316
$ws.Cells.Item(137, 9).Value = 981.7692
$ws.Cells.Item(137, 11).Value = 2945.3076
$ws.Cells.Item(137, 13).Value = -395.3076000000001

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45, 8).Value = 5758
$ws.Cells.Item(45, 9).Value = 5526.3335
$ws.Cells.Item(45, 11).Value = 5526.3335
$ws.Cells.Item(45, 13).Value = -5149.3335
$ws.Cells.Item(63, 8).Value = 2890
$ws.Cells.Item(63, 10).Value = 3800
$ws.Cells.Item(63, 12).Value = 3800
$ws.Cells.Item(63, 14).Value = -5172
$ws.Cells.Item(66, 8).Value = 2890
$ws.Cells.Item(66, 10).Value = 3800
$ws.Cells.Item(66, 12).Value = 19000
$ws.Cells.Item(66, 14).Value = -25864
$ws.Cells.Item(74, 8).Value = 3111.682
$ws.Cells.Item(74, 9).Value = 2049
$ws.Cells.Item(74, 10).Value = 4386.9
$ws.Cells.Item(74, 11).Value = 2049
$ws.Cells.Item(74, 12).Value = 4386.9
$ws.Cells.Item(74, 13).Value = -1175
$ws.Cells.Item(74, 14).Value = -6134.9
$ws.Cells.Item(77, 8).Value = 3111.682
$ws.Cells.Item(77, 9).Value = 2049
$ws.Cells.Item(77, 10).Value = 4386.9
$ws.Cells.Item(77, 11).Value = 10245
$ws.Cells.Item(77, 12).Value = 21934.5
$ws.Cells.Item(77, 13).Value = -5877
$ws.Cells.Item(77, 14).Value = -30670.5
$ws.Cells.Item(80, 8).Value = 18028.572
$ws.Cells.Item(80, 10).Value = 18028.572
$ws.Cells.Item(80, 12).Value = 18028.572
$ws.Cells.Item(80, 14).Value = -20024.572
$ws.Cells.Item(83, 8).Value = 18028.572
$ws.Cells.Item(83, 10).Value = 18028.572
$ws.Cells.Item(83, 12).Value = 54085.716
$ws.Cells.Item(83, 14).Value = -64069.716
$ws.Cells.Item(124, 8).Value = 26342.125
$ws.Cells.Item(124, 10).Value = 26342.125
$ws.Cells.Item(124, 12).Value = 26342.125
$ws.Cells.Item(124, 14).Value = -36162.125
$ws.Cells.Item(132, 8).Value = 3468.8572
$ws.Cells.Item(132, 9).Value = 3611.818
$ws.Cells.Item(132, 10).Value = 2944.6667
$ws.Cells.Item(132, 11).Value = 10835.454
$ws.Cells.Item(132, 12).Value = 8834.000100000001
$ws.Cells.Item(132, 13).Value = -8305.454000000002
$ws.Cells.Item(132, 14).Value = -13894.0001

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(82, 8).Value = 17783.385
$ws.Cells.Item(82, 10).Value = 29858.428
$ws.Cells.Item(82, 12).Value = 29858.428
$ws.Cells.Item(82, 14).Value = -30624.428
$ws.Cells.Item(85, 8).Value = 17783.385
$ws.Cells.Item(85, 10).Value = 29858.428
$ws.Cells.Item(85, 12).Value = 29858.428
$ws.Cells.Item(85, 14).Value = -32510.428
$ws.Cells.Item(86, 8).Value = 66625.82000000001
$ws.Cells.Item(86, 9).Value = 123522.445
$ws.Cells.Item(86, 10).Value = 2617.125
$ws.Cells.Item(86, 11).Value = 123522.445
$ws.Cells.Item(86, 12).Value = 2617.125
$ws.Cells.Item(86, 13).Value = -122399.445
$ws.Cells.Item(86, 14).Value = -4863.125
$ws.Cells.Item(89, 8).Value = 66625.82000000001
$ws.Cells.Item(89, 9).Value = 123522.445
$ws.Cells.Item(89, 10).Value = 2617.125
$ws.Cells.Item(89, 11).Value = 617612.2250000001
$ws.Cells.Item(89, 12).Value = 13085.625
$ws.Cells.Item(89, 13).Value = -611996.2250000001
$ws.Cells.Item(89, 14).Value = -24317.625
$ws.Cells.Item(134, 8).Value = 3107.5715
$ws.Cells.Item(134, 9).Value = 3090.818
$ws.Cells.Item(134, 10).Value = 3169
$ws.Cells.Item(134, 11).Value = 9272.454000000002
$ws.Cells.Item(134, 12).Value = 9507
$ws.Cells.Item(134, 13).Value = -6737.454000000002
$ws.Cells.Item(134, 14).Value = -14577

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 35729.4
$ws.Cells.Item(31, 9).Value = 41194.16
$ws.Cells.Item(31, 10).Value = 8405.6
$ws.Cells.Item(31, 11).Value = 41194.16
$ws.Cells.Item(31, 12).Value = 8405.6
$ws.Cells.Item(31, 13).Value = -40899.16
$ws.Cells.Item(31, 14).Value = -8995.6
$ws.Cells.Item(34, 8).Value = 35729.4
$ws.Cells.Item(34, 9).Value = 41194.16
$ws.Cells.Item(34, 10).Value = 8405.6
$ws.Cells.Item(34, 11).Value = 41194.16
$ws.Cells.Item(34, 12).Value = 8405.6
$ws.Cells.Item(34, 13).Value = -40992.16
$ws.Cells.Item(34, 14).Value = -8809.6
$ws.Cells.Item(58, 8).Value = 15938.111
$ws.Cells.Item(58, 9).Value = 1691.4667
$ws.Cells.Item(58, 10).Value = 87171.336
$ws.Cells.Item(58, 11).Value = 1691.4667
$ws.Cells.Item(58, 12).Value = 87171.336
$ws.Cells.Item(58, 13).Value = -1488.4667
$ws.Cells.Item(58, 14).Value = -87577.336
$ws.Cells.Item(74, 8).Value = 22649.8
$ws.Cells.Item(74, 10).Value = 22649.8
$ws.Cells.Item(74, 12).Value = 22649.8
$ws.Cells.Item(74, 14).Value = -24397.8
$ws.Cells.Item(77, 8).Value = 22649.8
$ws.Cells.Item(77, 10).Value = 22649.8
$ws.Cells.Item(77, 12).Value = 67949.39999999999
$ws.Cells.Item(77, 14).Value = -76685.39999999999
$ws.Cells.Item(99, 8).Value = 28448.5
$ws.Cells.Item(99, 9).Value = 4490
$ws.Cells.Item(99, 10).Value = 52407
$ws.Cells.Item(99, 11).Value = 4490
$ws.Cells.Item(99, 12).Value = 52407
$ws.Cells.Item(99, 13).Value = -2992
$ws.Cells.Item(99, 14).Value = -55403
$ws.Cells.Item(107, 8).Value = 652.9286
$ws.Cells.Item(107, 9).Value = 718.8570999999999
$ws.Cells.Item(107, 10).Value = 587
$ws.Cells.Item(107, 11).Value = 718.8570999999999
$ws.Cells.Item(107, 12).Value = 587
$ws.Cells.Item(107, 13).Value = 1201.1429
$ws.Cells.Item(107, 14).Value = -4427
$ws.Cells.Item(126, 8).Value = 28448.5
$ws.Cells.Item(126, 9).Value = 4490
$ws.Cells.Item(126, 10).Value = 52407
$ws.Cells.Item(126, 11).Value = 13470
$ws.Cells.Item(126, 12).Value = 157221
$ws.Cells.Item(126, 13).Value = -11000
$ws.Cells.Item(126, 14).Value = -162161
$ws.Cells.Item(136, 8).Value = 15938.111
$ws.Cells.Item(136, 9).Value = 1691.4667
$ws.Cells.Item(136, 10).Value = 87171.336
$ws.Cells.Item(136, 11).Value = 5074.4001
$ws.Cells.Item(136, 12).Value = 261514.008
$ws.Cells.Item(136, 13).Value = -2524.4001
$ws.Cells.Item(136, 14).Value = -266614.008
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(39, 8).Value = 3850
$ws.Cells.Item(39, 10).Value = 3850
$ws.Cells.Item(39, 12).Value = 11550
$ws.Cells.Item(39, 14).Value = -12138
$ws.Cells.Item(55, 8).Value = 9007.706
$ws.Cells.Item(55, 10).Value = 9533.1875
$ws.Cells.Item(55, 12).Value = 28599.5625
$ws.Cells.Item(55, 14).Value = -28953.5625
$ws.Cells.Item(131, 8).Value = 845.47
$ws.Cells.Item(131, 10).Value = 875.8817
$ws.Cells.Item(131, 12).Value = 2627.6451
$ws.Cells.Item(131, 14).Value = -12707.6451

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 63831.44
$ws.Cells.Item(70, 9).Value = 88036.28999999999
$ws.Cells.Item(70, 10).Value = 5739.8
$ws.Cells.Item(70, 11).Value = 88036.28999999999
$ws.Cells.Item(70, 12).Value = 5739.8
$ws.Cells.Item(70, 13).Value = -87766.28999999999
$ws.Cells.Item(70, 14).Value = -6279.8
$ws.Cells.Item(73, 8).Value = 63831.44
$ws.Cells.Item(73, 9).Value = 88036.28999999999
$ws.Cells.Item(73, 10).Value = 5739.8
$ws.Cells.Item(73, 11).Value = 88036.28999999999
$ws.Cells.Item(73, 12).Value = 5739.8
$ws.Cells.Item(73, 13).Value = -87100.28999999999
$ws.Cells.Item(73, 14).Value = -7611.8
$ws.Cells.Item(132, 8).Value = 2310.4119
$ws.Cells.Item(132, 9).Value = 1559.9231
$ws.Cells.Item(132, 11).Value = 4679.7693
$ws.Cells.Item(132, 13).Value = -2149.7693

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 112755.78
$ws.Cells.Item(16, 9).Value = 167133.5
$ws.Cells.Item(16, 10).Value = 4000.3333
$ws.Cells.Item(16, 11).Value = 167133.5
$ws.Cells.Item(16, 12).Value = 4000.3333
$ws.Cells.Item(16, 13).Value = -166963.5
$ws.Cells.Item(16, 14).Value = -4340.3333
$ws.Cells.Item(46, 8).Value = 724071.4399999999
$ws.Cells.Item(46, 10).Value = 1447400
$ws.Cells.Item(46, 12).Value = 1447400
$ws.Cells.Item(46, 14).Value = -1447776

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 10933
$ws.Cells.Item(132, 9).Value = 12119.8
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 36359.39999999999
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -33829.39999999999
$ws.Cells.Item(132, 14).Value = -20057
$ws.Cells.Item(136, 8).Value = 2038.4286
$ws.Cells.Item(136, 9).Value = 638
$ws.Cells.Item(136, 10).Value = 2598.6
$ws.Cells.Item(136, 11).Value = 1914
$ws.Cells.Item(136, 12).Value = 7795.799999999999
$ws.Cells.Item(136, 13).Value = 636
$ws.Cells.Item(136, 14).Value = -12895.8
